# BAB V.docx conclusion/suggestion wording fixes + _GoBack bookmark relocation.

$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------------
# "Kesimpulan" paragraph: drop the stray _GoBack bookmark that used to sit
# between ", maka dapat " and "disimpulkan bahwa" so the sentence reads as
# one continuous run of text again.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rngConclusion = $d.Content
$rngConclusion.Find.Execute(", maka dapat disimpulkan bahwa", $false, $false, $false, $false, $false, $true, 1, $false, ", maka dapat disimpulkan bahwa", 2) | Out-Null

# --- Change 2 --------------------------------------------------------------
# "Saran" paragraph: reword the suggestion sentence.
$rngSaran = $d.Content
$rngSaran.Find.Execute("lanjut untuk aplikasi ini agar dapat berjalan dengan sempurna dengan fungsi yang lebih baik adalah sebagai berikut:", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "lanjut untuk aplikasi ini agar dapat berjalan lebih baik lagi adalah sebagai berikut:", 2) | Out-Null

# --- Change 3 --------------------------------------------------------------
# Re-add the _GoBack bookmark (as a collapsed/empty range) right after the
# second occurrence of " positif dan negatif", just before the closing
# period, in the "Merubah proses pelabelan ..." bullet.
$fullText = $d.Content.Text
$firstHit = $fullText.IndexOf(" positif dan negatif")
$secondHit = $fullText.IndexOf(" positif dan negatif", $firstHit + 1)
$bookmarkPos = $secondHit + (" positif dan negatif").Length

$rngBookmark = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $rngBookmark)
